$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Cell H10 currently holds a "blog" entry referencing ser: 163.
# Update it to reference ser: 167 instead (keeping the same type/width/height lines).
$ws.Range("H10").Value = "type: blog`nwidth: 2`nheight: 1`nser: 167"

# Move the active selection to C10 (matches the final saved selection state).
$ws.Range("C10").Select()
